# Update the carjacking-by-neighborhood-by-month workbook to add one more
# day of data (through 2022-12-22, previously through 2022-12-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet (tab) name reflects the new "through" date.
$ws.Name = "Through 2022-12-22"

# The header in B1 (shared string) reflects the new "through" date too.
$ws.Range("B1").Value = "December 2022 (through December 22)"

# Updated / newly-populated per-neighborhood counts for the current month
# (column B) and the handful of historical-month cells that changed.
$ws.Range("Z3").Value = 10
$ws.Range("AL3").Value = 4
$ws.Range("BJ3").Value = 8
$ws.Range("B4").Value = 2
$ws.Range("Z4").Value = 3
$ws.Range("N5").Value = 3
$ws.Range("BV5").Value = 7
$ws.Range("BV6").Value = 1
$ws.Range("B7").Value = 7
$ws.Range("BJ7").Value = 2
$ws.Range("N10").Value = 5
$ws.Range("CH10").Value = 1
$ws.Range("AX14").Value = 11
$ws.Range("B15").Value = 5
$ws.Range("Z15").Value = 3
$ws.Range("AX15").Value = 2
$ws.Range("AL20").Value = 6
$ws.Range("N21").Value = 2
$ws.Range("AL22").Value = 2
$ws.Range("AL24").Value = 1
$ws.Range("CH32").Value = 1
$ws.Range("N33").Value = 2
$ws.Range("B40").Value = 3
$ws.Range("AL41").Value = 2
$ws.Range("BJ42").Value = 3
$ws.Range("Z43").Value = 3
$ws.Range("AX51").Value = 2
$ws.Range("BV60").Value = 1
$ws.Range("N64").Value = 8
$ws.Range("BJ72").Value = 1
$ws.Range("Z83").Value = 1
$ws.Range("BJ91").Value = 1
$ws.Range("Z94").Value = 1
$ws.Range("BJ96").Value = 1
$ws.Range("B97").Value = 4
$ws.Range("N97").Value = 2
